# Reproduce the authored change:
#  - styles.xml: the two bold header/title fonts collapse into a single
#    bold + white font (used by both the dashboard title cell and the
#    blue header row), dropping the separate 14pt title font.
#  - Training Dashboard!H3: -47 -> -55
#  - Training Dashboard!I3: "08-Sep-2025" -> "16-Sep-2025" (kept as text)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used  = $ws.UsedRange
    $ncols = $used.Columns.Count

    # Title cell (row 1) - stays bold, loses its 14pt size, turns white
    $titleCell = $ws.Range("A1")
    $titleCell.Font.Size  = 11
    $titleCell.Font.Color = 16777215

    # Header row (row 2) - stays bold, turns white (on its dark blue fill)
    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $ncols))
    $headerRange.Font.Color = 16777215
}

$ws1 = $wb.Worksheets.Item(1)

# H3: -47 -> -55
$ws1.Range("H3").Value2 = -55

# I3: force text so "16-Sep-2025" isn't auto-converted to a date serial
$ws1.Range("I3").NumberFormat = "@"
$ws1.Range("I3").Value2 = "16-Sep-2025"
